$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 200.9
$ws.Range("I33").Value = 132.375
$ws.Range("J33").Value = 475
$ws.Range("K33").Value = 132.375
$ws.Range("L33").Value = 475
$ws.Range("M33").Value = 96.625
$ws.Range("N33").Value = -933
$ws.Range("H106").Value = 10867
$ws.Range("I106").Value = 12831.091
$ws.Range("K106").Value = 12831.091
$ws.Range("M106").Value = -12200.091
$ws.Range("H132").Value = 7095946.5
$ws.Range("I132").Value = 9528074
$ws.Range("J132").Value = 2242.5833
$ws.Range("K132").Value = 28584222
$ws.Range("L132").Value = 6727.749899999999
$ws.Range("M132").Value = -28581692
$ws.Range("N132").Value = -11787.7499
$ws.Range("H137").Value = 1316.1471
$ws.Range("I137").Value = 975.65216
$ws.Range("K137").Value = 2926.95648
$ws.Range("M137").Value = -376.9564799999998
$ws.Range("H138").Value = 1373.28
$ws.Range("I138").Value = 804.64703
$ws.Range("J138").Value = 1666.2122
$ws.Range("K138").Value = 2413.94109
$ws.Range("L138").Value = 4998.6366
$ws.Range("M138").Value = 2726.05891
$ws.Range("N138").Value = -15278.6366
$ws.Range("H140").Value = 33337.145
$ws.Range("J140").Value = 33337.145
$ws.Range("L140").Value = 33337.145
$ws.Range("N140").Value = -43697.145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4679.8
$ws.Range("I2").Value = 598.3333
$ws.Range("J2").Value = 15175
$ws.Range("K2").Value = 598.3333
$ws.Range("L2").Value = 15175
$ws.Range("M2").Value = -485.3333
$ws.Range("N2").Value = -15401
$ws.Range("H32").Value = 4694.1885
$ws.Range("I32").Value = 4386.196
$ws.Range("J32").Value = 6718.143
$ws.Range("K32").Value = 4386.196
$ws.Range("L32").Value = 6718.143
$ws.Range("M32").Value = -4099.196
$ws.Range("N32").Value = -7292.143
$ws.Range("H45").Value = 1097.3334
$ws.Range("I45").Value = 1013.55
$ws.Range("J45").Value = 1264.9
$ws.Range("K45").Value = 1013.55
$ws.Range("L45").Value = 1264.9
$ws.Range("M45").Value = -636.55
$ws.Range("N45").Value = -2018.9
$ws.Range("H61").Value = 1215.9048
$ws.Range("I61").Value = 1119.4117
$ws.Range("J61").Value = 1626
$ws.Range("K61").Value = 1119.4117
$ws.Range("L61").Value = 1626
$ws.Range("M61").Value = -907.4117000000001
$ws.Range("N61").Value = -2050
$ws.Range("H102").Value = 41667492
$ws.Range("I102").Value = 41667492
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 41667492
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -41665870
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 1855.9231
$ws.Range("I110").Value = 1330.5714
$ws.Range("K110").Value = 1330.5714
$ws.Range("M110").Value = 714.4286
$ws.Range("H116").Value = 4679.8
$ws.Range("I116").Value = 598.3333
$ws.Range("J116").Value = 15175
$ws.Range("K116").Value = 598.3333
$ws.Range("L116").Value = 15175
$ws.Range("M116").Value = 1695.6667
$ws.Range("N116").Value = -19763
$ws.Range("H122").Value = 777
$ws.Range("I122").Value = 781.2
$ws.Range("K122").Value = 2343.6
$ws.Range("M122").Value = 106.3999999999996
$ws.Range("H136").Value = 1215.9048
$ws.Range("I136").Value = 1119.4117
$ws.Range("J136").Value = 1626
$ws.Range("K136").Value = 3358.2351
$ws.Range("L136").Value = 4878
$ws.Range("M136").Value = -808.2351000000003
$ws.Range("N136").Value = -9978

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4679.8
$ws.Range("I3").Value = 598.3333
$ws.Range("J3").Value = 15175
$ws.Range("K3").Value = 598.3333
$ws.Range("L3").Value = 15175
$ws.Range("M3").Value = -484.3333
$ws.Range("N3").Value = -15403
$ws.Range("H107").Value = 1712.4
$ws.Range("I107").Value = 1363.3334
$ws.Range("J107").Value = 3108.6667
$ws.Range("K107").Value = 1363.3334
$ws.Range("L107").Value = 3108.6667
$ws.Range("M107").Value = 556.6666
$ws.Range("N107").Value = -6948.6667
$ws.Range("H134").Value = 4641.528
$ws.Range("I134").Value = 1126.9656
$ws.Range("J134").Value = 19201.857
$ws.Range("K134").Value = 3380.8968
$ws.Range("L134").Value = 57605.571
$ws.Range("M134").Value = -845.8968
$ws.Range("N134").Value = -62675.571

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1643.9524
$ws.Range("I31").Value = 1166.1111
$ws.Range("J31").Value = 2002.3334
$ws.Range("K31").Value = 1166.1111
$ws.Range("L31").Value = 2002.3334
$ws.Range("M31").Value = -871.1111000000001
$ws.Range("N31").Value = -2592.3334
$ws.Range("H34").Value = 1643.9524
$ws.Range("I34").Value = 1166.1111
$ws.Range("J34").Value = 2002.3334
$ws.Range("K34").Value = 1166.1111
$ws.Range("L34").Value = 2002.3334
$ws.Range("M34").Value = -964.1111000000001
$ws.Range("N34").Value = -2406.3334
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H86").Value = 6690408.5
$ws.Range("I86").Value = 22225290
$ws.Range("J86").Value = 32602
$ws.Range("K86").Value = 22225290
$ws.Range("L86").Value = 32602
$ws.Range("M86").Value = -22224167
$ws.Range("N86").Value = -34848
$ws.Range("H89").Value = 6690408.5
$ws.Range("I89").Value = 22225290
$ws.Range("J89").Value = 32602
$ws.Range("K89").Value = 111126450
$ws.Range("L89").Value = 163010
$ws.Range("M89").Value = -111120834
$ws.Range("N89").Value = -174242
$ws.Range("H105").Value = 767.41174
$ws.Range("I105").Value = 779.8182
$ws.Range("J105").Value = 744.6667
$ws.Range("K105").Value = 779.8182
$ws.Range("L105").Value = 744.6667
$ws.Range("M105").Value = 967.1818
$ws.Range("N105").Value = -4238.6667
$ws.Range("H107").Value = 503.7586
$ws.Range("I107").Value = 471.9375
$ws.Range("J107").Value = 542.9231
$ws.Range("K107").Value = 471.9375
$ws.Range("L107").Value = 542.9231
$ws.Range("M107").Value = 1448.0625
$ws.Range("N107").Value = -4382.9231
$ws.Range("H132").Value = 1312.4103
$ws.Range("I132").Value = 748.5806
$ws.Range("K132").Value = 2245.7418
$ws.Range("M132").Value = 284.2582000000002
$ws.Range("H134").Value = 650.9524
$ws.Range("I134").Value = 494.07693
$ws.Range("J134").Value = 905.875
$ws.Range("K134").Value = 1482.23079
$ws.Range("L134").Value = 2717.625
$ws.Range("M134").Value = 1052.76921
$ws.Range("N134").Value = -7787.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 487921.75
$ws.Range("I4").Value = 224929.75
$ws.Range("K4").Value = 674789.25
$ws.Range("M4").Value = -674677.25
$ws.Range("H103").Value = 9865.75
$ws.Range("I103").Value = 697.5
$ws.Range("J103").Value = 14449.875
$ws.Range("K103").Value = 2092.5
$ws.Range("L103").Value = 43349.625
$ws.Range("M103").Value = -1213.5
$ws.Range("N103").Value = -45107.625
$ws.Range("H131").Value = 1067.36
$ws.Range("J131").Value = 1115.8298
$ws.Range("L131").Value = 3347.4894
$ws.Range("N131").Value = -13427.4894
$ws.Range("H138").Value = 1692.619
$ws.Range("I138").Value = 1318.091
$ws.Range("J138").Value = 2104.6
$ws.Range("K138").Value = 3954.273
$ws.Range("L138").Value = 6313.799999999999
$ws.Range("M138").Value = 1185.727
$ws.Range("N138").Value = -16593.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1882.0303
$ws.Range("I132").Value = 1155.409
$ws.Range("J132").Value = 3335.2727
$ws.Range("K132").Value = 3466.227
$ws.Range("L132").Value = 10005.8181
$ws.Range("M132").Value = -936.2270000000003
$ws.Range("N132").Value = -15065.8181

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H40").Value = 4799.273
$ws.Range("I40").Value = 3079.4
$ws.Range("K40").Value = 3079.4
$ws.Range("M40").Value = -2943.4
$ws.Range("H136").Value = 1935.0667
$ws.Range("I136").Value = 1894.3572
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 5683.071599999999
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -3133.071599999999
$ws.Range("N136").Value = -12615

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 415.3
$ws.Range("I107").Value = 397.46155
$ws.Range("K107").Value = 1192.38465
$ws.Range("M107").Value = 727.61535
$ws.Range("H132").Value = 1812.6666
$ws.Range("I132").Value = 1453.7778
$ws.Range("K132").Value = 4361.3334
$ws.Range("M132").Value = -1831.3334
$ws.Range("H135").Value = 39905
$ws.Range("J135").Value = 39905
$ws.Range("L135").Value = 39905
$ws.Range("N135").Value = -50045
$ws.Range("H140").Value = 32162.9
$ws.Range("J140").Value = 32162.9
$ws.Range("L140").Value = 32162.9
$ws.Range("N140").Value = -42522.9
